$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 97. This shifts the existing rows 97-133
# down to 98-134, matching the diff (new dimension A1:R134).
$ws.Rows("97").Insert()

# Populate the newly inserted row 97 with the latest weekly observation.
$ws.Range("A97").Value = 4
$ws.Range("B97").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C97").Value = "Los Lagos"
$ws.Range("D97").Value = 44468
$ws.Range("E97").Value = 10
$ws.Range("F97").Value = 100112017
$ws.Range("G97").Value = "Apio"
$ws.Range("H97").Value = "Americana (o)"
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 15
$ws.Range("K97").Value = 12000
$ws.Range("L97").Value = 12000
$ws.Range("M97").Value = 12000
$ws.Range("N97").Value = "$/docena de matas"
$ws.Range("O97").Value = "Región de Coquimbo"
$ws.Range("P97").Value = 2000
$ws.Range("Q97").Value = 6
$ws.Range("R97").Value = "Hortaliza"
